$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of notes (row 12) about subnational weighted averaging
$ws.Range("A12").Value = [DateTime]"2021-01-08"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "Subnational coverage estimates (mean, lower, upper) were incorrectly combined using a population-weighted average."
$ws.Range("F12").Value = "This should be corrected"

$ws.Rows.Item(12).RowHeight = 45

$ws.Range("A13").Select()
